# Generate Report for Handback
#
# The localization-status report is regenerated: the entry for
# "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md" moved from "Ready for handoff"
# (with a stale-handback error) to "Handed back: in sync with en-US", with
# refreshed handoff/handback timestamps; the three rows on every sheet are
# re-sorted (c7dafe61 first, then ffff25e0..., then ffffffbbafeb98...).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": A1:G4
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws.Range("B2").Value = "e2e\c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws.Range("C2").Value = ".md"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "Handed back: in sync with en-US"
$ws.Range("G2").Value = "2016-09-01 09:23:26"

$ws.Range("A3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws.Range("B3").Value = "e2e\ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-09-01 09:20:15"

$ws.Range("A4").Value = "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
$ws.Range("B4").Value = "e2e\ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
$ws.Range("C4").Value = ".md"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "Handed back: in sync with en-US"
$ws.Range("G4").Value = "2016-09-01 09:20:15"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = "e2e\c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
    } elseif ($addr -eq '$B$3') {
        $h.TextToDisplay = "e2e\ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
    } elseif ($addr -eq '$B$4') {
        $h.TextToDisplay = "e2e\ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn": A1:P4
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.1d934907a74ac1423164f5eb0eb4fb60ad0e471f.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-01 09:23:21"
$ws.Range("I2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws.Range("J2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.1d934907a74ac1423164f5eb0eb4fb60ad0e471f.zh-cn.xlf"
$ws.Range("K2").Value = "2016-09-01 09:23:39"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "'True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "'False"
$ws.Range("P2").Value = ""

$ws.Range("A3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-01 09:19:58"
$ws.Range("I3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws.Range("J3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws.Range("K3").Value = "2016-09-01 09:20:36"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = ""

$ws.Range("A4").Value = "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "e2e"
$ws.Range("E4").Value = "ht"
$ws.Range("F4").Value = "'True"
$ws.Range("G4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws.Range("H4").Value = "2016-09-01 09:19:58"
$ws.Range("I4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws.Range("J4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws.Range("K4").Value = "2016-09-01 09:20:36"
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = "'True"
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = "'False"
$ws.Range("P4").Value = ""

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
    } elseif ($addr -eq '$A$4') {
        $h.TextToDisplay = "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
    } elseif ($addr -eq '$I$4') {
        $h.TextToDisplay = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
    }
}

$ws.Columns.Item(16).ColumnWidth = 13.7470528738839

# ---------------------------------------------------------------------
# Sheet "de-de": A1:P4
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.1d934907a74ac1423164f5eb0eb4fb60ad0e471f.de-de.xlf"
$ws.Range("H2").Value = "2016-09-01 09:23:26"
$ws.Range("I2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws.Range("J2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.1d934907a74ac1423164f5eb0eb4fb60ad0e471f.de-de.xlf"
$ws.Range("K2").Value = "2016-09-01 09:23:46"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "'True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "'False"
$ws.Range("P2").Value = ""

$ws.Range("A3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws.Range("H3").Value = "2016-09-01 09:20:15"
$ws.Range("I3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws.Range("J3").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws.Range("K3").Value = "2016-09-01 09:20:43"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = ""

$ws.Range("A4").Value = "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "e2e"
$ws.Range("E4").Value = "ht"
$ws.Range("F4").Value = "'True"
$ws.Range("G4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws.Range("H4").Value = "2016-09-01 09:20:15"
$ws.Range("I4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws.Range("J4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws.Range("K4").Value = "2016-09-01 09:20:43"
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = "'True"
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = "'False"
$ws.Range("P4").Value = ""

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
    } elseif ($addr -eq '$A$4') {
        $h.TextToDisplay = "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
    } elseif ($addr -eq '$I$4') {
        $h.TextToDisplay = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
    }
}

$ws.Columns.Item(16).ColumnWidth = 13.7470528738839
